$d = $word.ActiveDocument

# 1. Merge "LRFD 2005 and later" + " with PT" into a single run/text.
$d.Content.Find.Execute(
    "LRFD 2005 and later with PT", $true, $false, $false, $false, $false,
    $true, 1, $false, "LRFD 2005 and later with PT", 2) | Out-Null

# 2. Merge "Before LRFD 2005" + " with PT" into a single run/text.
$d.Content.Find.Execute(
    "Before LRFD 2005 with PT", $true, $false, $false, $false, $false,
    $true, 1, $false, "Before LRFD 2005 with PT", 2) | Out-Null

# 3. Remove the trailing "_GoBack" bookmark (it disappears entirely in the
#    target revision).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 4. Insert the new trailing content: an empty paragraph, a "9th Edition"
#    paragraph, and a new equation (Aps*fps > As*fy) paragraph. Inserting at
#    the (collapsed) start of the final paragraph pushes the first N-1
#    paragraphs in front of it as new siblings, while the content of the
#    last inserted paragraph is folded into the existing final paragraph
#    (which keeps its identity, now bookmark-free).
$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">
<w:body>
<w:p/>
<w:p>
<w:r><w:t>9</w:t></w:r>
<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>
<w:r><w:t xml:space="preserve"> Edition</w:t></w:r>
</w:p>
<w:p>
<m:oMathPara>
<m:oMath>
<m:sSub>
<m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr>
<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>A</m:t></m:r></m:e>
<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>ps</m:t></m:r></m:sub>
</m:sSub>
<m:sSub>
<m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr>
<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>f</m:t></m:r></m:e>
<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>ps</m:t></m:r></m:sub>
</m:sSub>
<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>&gt;</m:t></m:r>
<m:sSub>
<m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr>
<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>A</m:t></m:r></m:e>
<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>s</m:t></m:r></m:sub>
</m:sSub>
<m:sSub>
<m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr>
<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>f</m:t></m:r></m:e>
<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>y</m:t></m:r></m:sub>
</m:sSub>
</m:oMath>
</m:oMathPara>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$insertPoint.InsertXML($xml)
